$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.710.76"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "3.423.37"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'572.67"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.424.96"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "4.013.79"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'0.0000188"
$ws.Range("E15").Value = "  -4.21%  "
$ws.Range("D16").Value = "'27.68"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("D17").Value = "64.700.35"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "3.385.79"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("D20").Value = "'13.79"
$ws.Range("E20").Value = "  -3.55%  "
$ws.Range("D21").Value = "'380.54"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'72.13"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("E26").Value = "  -4.90%  "
$ws.Range("D27").Value = "'10.00"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").Value = "'0.178"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").Value = "'0.990"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("D32").Value = "'2.01"
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("D33").Value = "'23.20"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").Value = "'7.07"
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").Value = "'160.40"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").Value = "'1.91"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").Value = "2.908.93"
$ws.Range("E38").Value = "  -5.62%  "
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("D41").Value = "'26.38"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("D42").Value = "'4.57"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").Value = "'43.00"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "'0.771"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").Value = "'25.71"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "'317.13"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").Value = "'2.23"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("E49").Value = "  -5.17%  "
$ws.Range("D51").Value = "'6.52"
$ws.Range("E51").Value = "  -2.84%  "
